# Generate Report for Handback
# This localization-status report is regenerated: the handback for
# zh-cn/de-de is now in sync with en-US, the "Latest Handback DateTime"
# values are refreshed, and the previously-recorded "version is not the
# latest" error is cleared now that the issue is resolved.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: Status summary columns for each locale ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("K2").Value = "2016-09-02 18:55:36"
$wsZhCn.Range("P2").Value = ""
$wsZhCn.Range("P2").Style = "Normal"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("K2").Value = "2016-09-02 18:55:45"
$wsDeDe.Range("P2").Value = ""
$wsDeDe.Range("P2").Style = "Normal"

# Re-fit the columns whose contents changed length (Status and Error
# Detail) so the report reads cleanly.
$wsOverview.Columns.Item(5).AutoFit()
$wsOverview.Columns.Item(6).AutoFit()

$wsZhCn.Columns.Item(3).AutoFit()
$wsZhCn.Columns.Item(16).AutoFit()

$wsDeDe.Columns.Item(3).AutoFit()
$wsDeDe.Columns.Item(16).AutoFit()
